$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("erosion")
$ws.Activate()

# Update C2 value from 0.1 to 0 (resolve single matrix issue)
$ws.Range("C2").Value = 0

# Move active selection to C3 to match the saved workbook state
$ws.Range("C3").Select()
